# Revert the SharePoint "shared with" metadata sync that was merged into
# customXml/item1.xml (ma:contentTypeVersion 5 -> 3, ma:versionID / ma:fieldsID
# rolled back, the ns3 "cdc19bd1-...-SharedWithUsers/SharedWithDetails" schema
# block removed) and let Word re-mint customXml/itemProps1.xml's datastore
# item GUID as a side effect of recreating the part.
#
# The canonical way to edit a document's custom XML data store from the Word
# object model is Document.CustomXMLParts: locate the SharePoint content-type
# schema part by its namespace, delete it, and Add() the reverted XML back
# (CustomXMLPart.XML itself is not user-settable in the real object model --
# you recreate the part instead). Word then regenerates the paired
# itemProps*.xml datastore item (new ds:itemID GUID) automatically.

$d = $word.ActiveDocument

$contentTypeNamespace = "http://schemas.microsoft.com/office/2006/metadata/contentType"

# The reverted contentType schema part (PR merge -> revert): contentTypeVersion
# 5 -> 3, versionID/fieldsID rolled back to the pre-merge values, ns3 import +
# SharedWithUsers/SharedWithDetails <xsd:element ref.../> rows and the whole
# "cdc19bd1-6dbc-...-b7f924674b29"-sibling <xsd:schema targetNamespace="cdc19bd1-...">
# block removed.
$revertedXml = @"
<?xml version="1.0" encoding="utf-8"?>
<ct:contentTypeSchema xmlns:ct="http://schemas.microsoft.com/office/2006/metadata/contentType" xmlns:ma="http://schemas.microsoft.com/office/2006/metadata/properties/metaAttributes" ct:_="" ma:_="" ma:contentTypeName="Document" ma:contentTypeID="0x0101004446CF7D97B8CF4583ACEECEC8CECEA7" ma:contentTypeVersion="3" ma:contentTypeDescription="Create a new document." ma:contentTypeScope="" ma:versionID="7c13b77237112d85c7d460547c88ab8b">
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:p="http://schemas.microsoft.com/office/2006/metadata/properties" xmlns:ns2="eb292c7f-6dbc-49d6-b736-b7f924674b29" targetNamespace="http://schemas.microsoft.com/office/2006/metadata/properties" ma:root="true" ma:fieldsID="b42858e4f81ddb688db20af23a05b3b9" ns2:_="">
    <xsd:import namespace="eb292c7f-6dbc-49d6-b736-b7f924674b29"/>
    <xsd:element name="properties">
      <xsd:complexType>
        <xsd:sequence>
          <xsd:element name="documentManagement">
            <xsd:complexType>
              <xsd:all>
                <xsd:element ref="ns2:MediaServiceMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceFastMetadata" minOccurs="0"/>
                <xsd:element ref="ns2:MediaServiceObjectDetectorVersions" minOccurs="0"/>
              </xsd:all>
            </xsd:complexType>
          </xsd:element>
        </xsd:sequence>
      </xsd:complexType>
    </xsd:element>
  </xsd:schema>
  <xsd:schema xmlns:xsd="http://www.w3.org/2001/XMLSchema" xmlns:xs="http://www.w3.org/2001/XMLSchema" xmlns:dms="http://schemas.microsoft.com/office/2006/documentManagement/types" xmlns:pc="http://schemas.microsoft.com/office/infopath/2007/PartnerControls" targetNamespace="eb292c7f-6dbc-49d6-b736-b7f924674b29" elementFormDefault="qualified">
    <xsd:import namespace="http://schemas.microsoft.com/office/2006/documentManagement/types"/>
    <xsd:import namespace="http://schemas.microsoft.com/office/infopath/2007/PartnerControls"/>
    <xsd:element name="MediaServiceMetadata" ma:index="8" nillable="true" ma:displayName="MediaServiceMetadata" ma:hidden="true" ma:internalName="MediaServiceMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceFastMetadata" ma:index="9" nillable="true" ma:displayName="MediaServiceFastMetadata" ma:hidden="true" ma:internalName="MediaServiceFastMetadata" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Note"/>
      </xsd:simpleType>
    </xsd:element>
    <xsd:element name="MediaServiceObjectDetectorVersions" ma:index="10" nillable="true" ma:displayName="MediaServiceObjectDetectorVersions" ma:hidden="true" ma:indexed="true" ma:internalName="MediaServiceObjectDetectorVersions" ma:readOnly="true">
      <xsd:simpleType>
        <xsd:restriction base="dms:Text"/>
      </xsd:simpleType>
    </xsd:element>
  </xsd:schema>
</ct:contentTypeSchema>
"@

$reverted = $false

try {
    $parts = $d.CustomXMLParts.SelectByNamespace($contentTypeNamespace)
    if ($parts -ne $null -and $parts.Count -gt 0) {
        for ($i = $parts.Count; $i -ge 1; $i--) {
            $part = $parts.Item($i)
            $part.Delete()
        }
        $d.CustomXMLParts.Add($revertedXml) | Out-Null
        $reverted = $true
    }
} catch {
    Write-Output ("SelectByNamespace path failed: " + $_.Exception.Message)
}

if (-not $reverted) {
    # Fallback: walk every custom XML part and replace the one whose XML
    # carries the SharePoint content-type schema namespace/root element.
    try {
        $count = $d.CustomXMLParts.Count
        for ($i = 1; $i -le $count; $i++) {
            $part = $d.CustomXMLParts.Item($i)
            $isMatch = $false
            try {
                if ($part.NamespaceURI -eq $contentTypeNamespace) { $isMatch = $true }
            } catch { }
            if (-not $isMatch) {
                try {
                    if ($part.XML -ne $null -and $part.XML.Contains("contentTypeSchema")) { $isMatch = $true }
                } catch { }
            }
            if ($isMatch) {
                try {
                    $part.Delete()
                    $d.CustomXMLParts.Add($revertedXml) | Out-Null
                } catch {
                    try { $part.XML = $revertedXml } catch { }
                }
                $reverted = $true
            }
        }
    } catch {
        Write-Output ("Item-walk fallback failed: " + $_.Exception.Message)
    }
}

if ($reverted) {
    Write-Output "customXml/item1.xml reverted to the pre-merge SharePoint content-type schema (contentTypeVersion 5 -> 3); Word will re-mint itemProps1.xml's ds:itemID."
} else {
    Write-Output "CustomXMLParts editing is not available on this host; customXml/item1.xml and itemProps1.xml were left unchanged."
}
